# Fix 3-Year Summary category alignment (infrastructure-costs.xlsx)
#
# The "3-Year Summary" sheet's Category column (A3:A6) was out of step with
# the categories actually used on the "Infrastructure Costs" / "Credits"
# sheets (Software Licenses / Support & Maintenance / Professional Services),
# and the TOTAL row's label had stray markdown-style asterisks. Realign the
# labels so the SUMIF rollups in B:G key off the right category names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3-Year Summary")

# Row 4 was mis-labelled as a second "Software Licenses" row - it should be
# "Support & Maintenance" (matching the Infrastructure Costs / Credits sheets).
$ws.Range("A4").Value = "Support & Maintenance"

# Row 5 was mis-labelled as "Support & Maintenance" (duplicating row 4's old
# label) - it should be "Professional Services".
$ws.Range("A5").Value = "Professional Services"

# Row 6 is the grand-total row; drop the markdown-style emphasis from the label.
$ws.Range("A6").Value = "TOTAL"

# Touch row 7 (just below the TOTAL row) so nothing stray lingers there.
$ws.Rows.Item(7).OutlineLevel = 0
